{"js": "const replacements = [\n  [\"718\u00d75=\", \"919\u00d74=\"],\n  [\"519\u00d77=\", \"375\u00d73=\"],\n  [\"618\u00d75=\", \"222\u00d79=\"],\n  [\"137\u00d75=\", \"794\u00d76=\"],\n  [\"307\u00d77=\", \"368\u00d72=\"],\n  [\"676\u00d75=\", \"230\u00d73=\"],\n  [\"689\u00d73=\", \"601\u00d76=\"],\n  [\"897\u00d78=\", \"641\u00d76=\"],\n  [\"157\u00d75=\", \"561\u00d76=\"],\n  [\"901\u00d78=\", \"757\u00d78=\"],\n  [\"792\u00d79=\", \"463\u00d77=\"],\n  [\"809\u00d77=\", \"255\u00d75=\"],\n  [\"286\u00d77=\", \"603\u00d72=\"],\n  [\"495\u00d76=\", \"797\u00d76=\"],\n  [\"605\u00d76=\", \"689\u00d76=\"],\n  [\"858\u00d72=\", \"415\u00d76=\"],\n  [\"232\u00d77=\", \"257\u00d73=\"],\n  [\"455\u00d78=\", \"279\u00d72=\"],\n  [\"946\u00d78=\", \"889\u00d72=\"],\n  [\"342\u00d72=\", \"137\u00d74=\"],\n  [\"251\u00d76=\", \"943\u00d78=\"],\n  [\"679\u00d72=\", \"851\u00d78=\"],\n  [\"506\u00d78=\", \"568\u00d73=\"],\n  [\"684\u00d73=\", \"886\u00d75=\"],\n  [\"218\u00d73=\", \"914\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = '718\u00d75='; Replace = '919\u00d74=' },\n    @{ Find = '519\u00d77='; Replace = '375\u00d73=' },\n    @{ Find = '618\u00d75='; Replace = '222\u00d79=' },\n    @{ Find = '137\u00d75='; Replace = '794\u00d76=' },\n    @{ Find = '307\u00d77='; Replace = '368\u00d72=' },\n    @{ Find = '676\u00d75='; Replace = '230\u00d73=' },\n    @{ Find = '689\u00d73='; Replace = '601\u00d76=' },\n    @{ Find = '897\u00d78='; Replace = '641\u00d76=' },\n    @{ Find = '157\u00d75='; Replace = '561\u00d76=' },\n    @{ Find = '901\u00d78='; Replace = '757\u00d78=' },\n    @{ Find = '792\u00d79='; Replace = '463\u00d77=' },\n    @{ Find = '809\u00d77='; Replace = '255\u00d75=' },\n    @{ Find = '286\u00d77='; Replace = '603\u00d72=' },\n    @{ Find = '495\u00d76='; Replace = '797\u00d76=' },\n    @{ Find = '605\u00d76='; Replace = '689\u00d76=' },\n    @{ Find = '858\u00d72='; Replace = '415\u00d76=' },\n    @{ Find = '232\u00d77='; Replace = '257\u00d73=' },\n    @{ Find = '455\u00d78='; Replace = '279\u00d72=' },\n    @{ Find = '946\u00d78='; Replace = '889\u00d72=' },\n    @{ Find = '342\u00d72='; Replace = '137\u00d74=' },\n    @{ Find = '251\u00d76='; Replace = '943\u00d78=' },\n    @{ Find = '679\u00d72='; Replace = '851\u00d78=' },\n    @{ Find = '506\u00d78='; Replace = '568\u00d73=' },\n    @{ Find = '684\u00d73='; Replace = '886\u00d75=' },\n    @{ Find = '218\u00d73='; Replace = '914\u00d75=' }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $r.Find\n    $range.Find.Replacement.Text = $r.Replace\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $found = $range.Find.Execute($r.Find, $true, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $($r.Find)\"\n    }\n}\n"}
